$wb = $excel.ActiveWorkbook

# Rename the worksheets (tabs renamed to include the new numeric suffixes)
$newNames = @("summ24980915", "summ30158412", "summ35980124", "summ41710845", "summ46912963", "summ52665432", "summ58487474", "summ04192391", "summ09901147")

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]

    # Update the "CarAvailable" parameter label (row 29, column A) on every sheet
    if ($ws.Range("A29").Value2 -eq "CarAvailable") {
        $ws.Range("A29").Value = "CarOwnershipHH"
    }
}
